$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '43.131.53'
$ws.Range("E2").Value = '  +4.55%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '2.241.32'
$ws.Range("E3").Value = '  +2.85%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.22%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.76%  '

$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.620'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.86%  '

$ws.Range("B7").Value = 'Solana'
$ws.Range("C7").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '75.91'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +7.79%  '

$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.615'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.96%  '

$ws.Range("B10").Value = 'Avalanche'
$ws.Range("C10").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.14'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.22%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0935'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.50%  '

$ws.Range("B12").Value = 'OKB'
$ws.Range("C12").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.42'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.47%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.98'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.97%  '

$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.102'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.13%  '

$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.558.56'
$ws.Range("E15").Value = '  +2.26%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.71'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.47%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.244.36'
$ws.Range("E17").Value = '  +2.81%  '

$ws.Range("B18").Value = 'Polygon'
$ws.Range("C18").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.815'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.49%  '

$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '43.040.98'
$ws.Range("E19").Value = '  +4.88%  '

$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000105'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.47%  '

$ws.Range("B21").Value = 'Litecoin'
$ws.Range("C21").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.86%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.64%  '

$ws.Range("B23").Value = 'InternetComputer(DFINITY)'
$ws.Range("C23").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.44'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.45%  '

$ws.Range("B24").Value = 'ImmutableX'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.21'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +13.25%  '

$ws.Range("B25").Value = 'BitcoinCash'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '230.78'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.12%  '

$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.20%  '

$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.97'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.87%  '

$ws.Range("B28").Value = 'WEMIXToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.37'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.15%  '

$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.25'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.73%  '

$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.72'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +21.02%  '

$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.96'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.08%  '

$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.13'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.85%  '

$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.35'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.67%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0794'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.21%  '

$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.38'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.76%  '

$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.123'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.33%  '

$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.112'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.83%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.37'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.76%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0333'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +16.88%  '

$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '13.16'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.32%  '

$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.14'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.80%  '

$ws.Range("B42").Value = 'THORChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.59'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.96%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.200'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.62%  '

$ws.Range("B44").Value = 'MultiversX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '60.30'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.12%  '

$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '105.37'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.11%  '

$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.63'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.54%  '

$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0993'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.10%  '

$ws.Range("B48").Value = 'WOONetwork'
$ws.Range("C48").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.443'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +20.20%  '

$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.11'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.36%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.31'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.18%  '

$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.15'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.41%  '
